# Round the metric values shown in the CNV_metrics table to 3 decimal
# places, matching the commit "updating examples output with rounding
# implemented".
#
# The table (Tables.Item(1)) has 5 columns:
#   1: Metric name   2: Sample value   3: Sample QC   4: Control value   5: Control QC
# and the following rows contain the numeric metrics that need rounding:
#   Row 2: Percent diff (2 Mbp window)  -> Sample col2, Control col4
#   Row 3: Percent diff (6 Mbp window)  -> Sample col2, Control col4
#   Row 4: Correlation with label density -> Control col4 only (Sample col2 already rounded)
#   Row 5: Wave template correlation    -> Sample col2, Control col4

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(2, 2).Range.Text = "0.800"
$t.Cell(2, 4).Range.Text = "4.728"

$t.Cell(3, 2).Range.Text = "11.180"
$t.Cell(3, 4).Range.Text = "12.245"

$t.Cell(4, 4).Range.Text = "0.070"

$t.Cell(5, 2).Range.Text = "0.000"
$t.Cell(5, 4).Range.Text = "0.000"

Write-Host "Updated metric cells with rounded values."
